# Add the GitHub repo link to column F (rows 2-43) and update the
# worksheet's current selection, matching the "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill F2:F43 with the new github link (single bulk write so the row
# "spans" bookkeeping matches a normal Excel fill/paste operation).
$ws.Range("F2:F43").Value = "https://github.com/rolls65/MySQL_PRACTICES"

# Restore the sheet's selection/scroll state to the cell recorded in the
# saved workbook (G2) - this also updates the <selection> element.
[void]$ws.Range("G2").Select()
